$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "'0420194406718"
$ws.Range("B4").Value = "'0420194406719"
$ws.Range("B2").Value = "'1220194200667  "

$ws.Range("C8").Select()
